$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target C:H values for rows 2-21 after the edit (4 new rows inserted at top,
# the last 5 old rows dropped, net one row shorter: A1:H22 -> A1:H21)
$data = @(
    @(-1.03581714630127,1.381664276123047,-0.1787742376327515,0.0007635815418325,0.01328631862998,0.0305432621389627),
    @(-1.100839495658875,1.406517148017883,-0.2175595723092557,0.0154243474826216,-0.0059559359215199,0.0029016099870204),
    @(-1.232075214385986,1.368059515953064,-0.2088889628648757,0.0140499006956815,0.0103847095742821,0.0612392425537109),
    @(-1.198645412921905,1.359035015106201,-0.2234921492636203,0.0195476878434419,0.0311541277915239,0.0459676086902618),
    @(-1.162086248397827,1.317891120910645,0.1167446374893195,0.0041233403608202,-0.0007635815418325,-0.066737025976181),
    @(-1.455766379833223,1.126043200492858,0.7863338142633457,-0.0314595587551593,0.2884811162948608,0.06536258012056349),
    @(-1.653480172157284,0.8844107389450073,1.495913922786714,-0.1331686228513717,0.493731826543808,-0.1162171140313148),
    @(-1.116286456584938,0.6627160906791623,2.260028153657919,-0.09651670604944219,0.7289149761199951,-0.0355829000473022),
    @(-3.20666265487671,-0.9849638938903851,3.828832626342773,-0.2052507251501083,0.5688682794570923,0.1585195362567901),
    @(-3.492754817008973,-1.960709273815156,3.70500636100769,-0.18539759516716,-0.1214094683527946,-0.6217080950737),
    @(-2.38550305366516,0.4144415855407741,-0.4261573851108604,0.8633053302764893,-0.9155342578887939,0.0097738439217209),
    @(-2.568133831024171,3.730020523071294,-0.8093817904591537,-2.978273391723633,1.564731359481812,-3.266449213027954),
    @(-1.976609468460079,6.545797109603887,2.459241539239894,1.346346974372864,1.957823157310486,1.074512004852295),
    @(-1.197644114494325,5.058232277631741,2.450196892023069,0.2492330223321914,0.9819658994674684,-0.9859365224838256),
    @(-1.255056142807007,0.5559926331043243,-1.602012172341346,0.0514653958380222,-0.3197879493236542,0.1852448880672454),
    @(-0.5990372896194439,1.264416024088864,-0.3774302378296863,0.1069014146924018,-0.1577559560537338,0.1151480972766876),
    @(-0.1129188537597673,2.17881894111634,-0.6807380914687996,0.7513642311096191,-1.080926060676575,-0.1253800988197326),
    @(-0.4983874559402485,3.747065991163262,1.904177859425558,-0.155312493443489,0.2770273983478546,0.052381694316864),
    @(-0.2239453792572,0.9145344123244095,-1.341536760330224,-0.0971275717020034,-0.6624833345413208,0.3292563557624817),
    @(0.0572257041931156,1.018226306885483,0.3575173318386149,-0.0768163055181503,-0.087353728711605,0.0009162978967650999)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Count; $j++) {
        $col = $j + 3  # C=3
        $ws.Cells.Item($r, $col).Value = $row[$j]
    }
}

# Remove the now-obsolete last data row (row 22), shrinking the sheet by one row
$ws.Rows.Item(22).Delete()
